$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: Wesley
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Wesley"
$ws.Range("D13").Value = '"158760334335672320"'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "11/12/2020"

# Row 14: Sydney
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Sydney"
$ws.Range("C14").NumberFormat = "mm-dd-yy"
$ws.Range("C14").Value = 43962
$ws.Range("D14").Value = '"814697817989185537"'

$ws.Range("Q24").Select()
